$wb = $excel.ActiveWorkbook

# --- Sheet "Tex, Flags<0>" (sheet1): add row 19 "Separate W Clip" ---
$ws1 = $wb.Worksheets.Item("Tex, Flags<0>")

$ws1.Range("A19").Value2 = "Separate W Clip"
$ws1.Range("B19").Value2 = 107261
$ws1.Range("C19").Value2 = 98039
$ws1.Range("D19").Value2 = 114
$ws1.Range("E19").Formula = "=(D19/D`$2)-1"
$ws1.Range("F19").Formula = "=(D19/D18)-1"
$ws1.Range("H19").Value2 = 15796

# K26/K27/K28 benchmark timing values
$ws1.Range("K26").Value2 = 8.7449999999999992
$ws1.Range("K27").Value2 = 8.7420000000000009
$ws1.Range("K28").Value2 = 8.74

# --- Sheet "Flat, Flags<0>" (sheet2): add row 13 "Separate W Clip" ---
$ws2 = $wb.Worksheets.Item("Flat, Flags<0>")

$ws2.Range("A13").Value2 = "Separate W Clip"
$ws2.Range("B13").Value2 = 267236
$ws2.Range("C13").Value2 = 265251
$ws2.Range("D13").Value2 = 557
$ws2.Range("E13").Formula = "=(D13/D`$2)-1"
$ws2.Range("F13").Formula = "=(D13/D12)-1"
$ws2.Range("H13").Value2 = 15796

# --- Selections: set per-sheet active cell, keep sheet1 as the active tab ---
[void]$ws2.Range("A14").Select()
[void]$ws1.Activate()
[void]$ws1.Range("A19").Select()
